$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = '28.673.73'
$ws.Range("E2").Value = '  +2.32%  '

$ws.Range("D3").Value = '1.800.37'
$ws.Range("E3").Value = '  -0.21%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.15%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '313.56'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.73%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.000'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.19%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5398'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -1.58%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3782'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.39%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07531'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.80%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '42.58'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.26%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.117'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.51%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.000'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.14%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '20.98'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.89%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.177'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.88%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.416'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.64%  '

$ws.Range("D16").Value = '1.797.05'
$ws.Range("E16").Value = '  -0.23%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '90.65'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.59%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001065'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.37%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06440'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.95%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.000'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.21%  '

$ws.Range("E21").Value = '  +0.55%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.939'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.92%  '

$ws.Range("D23").Value = '28.670.65'
$ws.Range("E23").Value = '  +2.20%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.17'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.67%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.100'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.22%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '160.63'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.87%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.49'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.57%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.374'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.40%  '

$ws.Range("D29").Value = '2.005.30'
$ws.Range("E29").Value = '  -0.20%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '123.42'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.29%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.111'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.60%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.1041'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.36%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.668'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.55%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.700'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.82%  '

$ws.Range("E35").Value = '  +6.37%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.06497'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +7.33%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '8.904'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +3.13%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02320'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.84%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.046'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.40%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '11.34'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.61%  '

$ws.Range("B41").Value = 'TheSandbox'
$ws.Range("C41").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6252'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.85%  '

$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.209'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +4.92%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.9999'
$ws.Range("D43").Style = "Normal"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.394'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.23%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.36'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.58%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5879'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.82%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.668'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.06%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '126.14'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +3.31%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.962'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.85%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.157'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.71%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06883'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.37%  '
